$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08230999999999999
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("M2").Value = 118.0346986666667
$ws.Range("N2").Value = 354.104096
$ws.Range("O2").Value = 0.2666057129183408
$ws.Range("P2").Value = 0.2666057129183408
$ws.Range("Q2").Value = 3.238478682417777
$ws.Range("R2").Value = 29.14630814176
$ws.Range("S2").Value = 0.001963893679055197
$ws.Range("T2").Value = 0.001963893679055197
# Row 3
$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08230999999999999
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.4881754016778185
$ws.Range("P3").Value = 0.4881754016778186
$ws.Range("Q3").Value = 5.92990155503
$ws.Range("R3").Value = 53.36911399527
$ws.Range("S3").Value = 0.003596039166343554
$ws.Range("T3").Value = 0.003596039166343554
# Row 4
$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08230999999999999
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("M4").Value = 45.876452
$ws.Range("N4").Value = 137.629356
$ws.Range("O4").Value = 0.1036214293744632
$ws.Range("P4").Value = 0.1036214293744632
$ws.Range("Q4").Value = 1.258696921373333
$ws.Range("R4").Value = 11.32827229236
$ws.Range("S4").Value = 0.0007633049867371129
$ws.Range("T4").Value = 0.000763304986737113
# Row 5
$ws.Range("G5").Value = 0.02743666666666666
$ws.Range("H5").Value = 0.08230999999999999
$ws.Range("I5").Value = 0.007366285056527356
$ws.Range("J5").Value = 0.007366285056527356
$ws.Range("M5").Value = 62.68962833333333
$ws.Range("N5").Value = 188.068885
$ws.Range("O5").Value = 0.1415974560293775
$ws.Range("P5").Value = 0.1415974560293775
$ws.Range("Q5").Value = 1.719994436038889
$ws.Range("R5").Value = 15.47994992435
$ws.Range("S5").Value = 0.001043047224391492
$ws.Range("T5").Value = 0.001043047224391493
# Row 6
$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("M6").Value = 118.0346986666667
$ws.Range("N6").Value = 354.104096
$ws.Range("O6").Value = 0.2666057129183408
$ws.Range("P6").Value = 0.2666057129183408
$ws.Range("Q6").Value = 397.5796985251946
$ws.Range("R6").Value = 3578.217286726751
$ws.Range("S6").Value = 0.2411021758745588
$ws.Range("T6").Value = 0.2411021758745588
# Row 7
$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("O7").Value = 0.4881754016778185
$ws.Range("P7").Value = 0.4881754016778186
$ws.Range("Q7").Value = 727.9987622993309
$ws.Range("R7").Value = 6551.988860693978
$ws.Range("S7").Value = 0.4414764794969318
$ws.Range("T7").Value = 0.4414764794969318
# Row 8
$ws.Range("G8").Value = 3.368329
$ws.Range("I8").Value = 0.9043398704228307
$ws.Range("J8").Value = 0.9043398704228307
$ws.Range("M8").Value = 45.876452
$ws.Range("N8").Value = 137.629356
$ws.Range("O8").Value = 0.1036214293744632
$ws.Range("P8").Value = 0.1036214293744632
$ws.Range("Q8").Value = 154.526983688708
$ws.Range("R8").Value = 1390.742853198372
$ws.Range("S8").Value = 0.09370899001353054
$ws.Range("T8").Value = 0.09370899001353056
# Row 9
$ws.Range("G9").Value = 3.368329
$ws.Range("I9").Value = 0.9043398704228307
$ws.Range("J9").Value = 0.9043398704228307
$ws.Range("M9").Value = 62.68962833333333
$ws.Range("N9").Value = 188.068885
$ws.Range("O9").Value = 0.1415974560293775
$ws.Range("P9").Value = 0.1415974560293775
$ws.Range("Q9").Value = 211.1592931143883
$ws.Range("R9").Value = 1900.433638029495
$ws.Range("S9").Value = 0.1280522250378097
$ws.Range("T9").Value = 0.1280522250378097
# Row 10
$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("M10").Value = 118.0346986666667
$ws.Range("N10").Value = 354.104096
$ws.Range("O10").Value = 0.2666057129183408
$ws.Range("P10").Value = 0.2666057129183408
$ws.Range("Q10").Value = 38.81708772801777
$ws.Range("R10").Value = 349.3537895521599
$ws.Range("S10").Value = 0.0235396433647269
$ws.Range("T10").Value = 0.0235396433647269
# Row 11
$ws.Range("G11").Value = 0.3288616666666667
$ws.Range("H11").Value = 0.9865849999999999
$ws.Range("I11").Value = 0.08829384452064198
$ws.Range("J11").Value = 0.08829384452064198
$ws.Range("O11").Value = 0.4881754016778185
$ws.Range("P11").Value = 0.4881754016778186
$ws.Range("Q11").Value = 71.077049273105
$ws.Range("R11").Value = 639.6934434579449
$ws.Range("S11").Value = 0.04310288301454326
$ws.Range("T11").Value = 0.04310288301454326
# Row 12
$ws.Range("G12").Value = 0.3288616666666667
$ws.Range("H12").Value = 0.9865849999999999
$ws.Range("I12").Value = 0.08829384452064198
$ws.Range("J12").Value = 0.08829384452064198
$ws.Range("M12").Value = 45.876452
$ws.Range("N12").Value = 137.629356
$ws.Range("O12").Value = 0.1036214293744632
$ws.Range("P12").Value = 0.1036214293744632
$ws.Range("Q12").Value = 15.08700646547333
$ws.Range("R12").Value = 135.78305818926
$ws.Range("S12").Value = 0.009149134374195537
$ws.Range("T12").Value = 0.009149134374195537
# Row 13
$ws.Range("G13").Value = 0.3288616666666667
$ws.Range("H13").Value = 0.9865849999999999
$ws.Range("I13").Value = 0.08829384452064198
$ws.Range("J13").Value = 0.08829384452064198
$ws.Range("M13").Value = 62.68962833333333
$ws.Range("N13").Value = 188.068885
$ws.Range("O13").Value = 0.1415974560293775
$ws.Range("P13").Value = 0.1415974560293775
$ws.Range("Q13").Value = 20.61621565641389
$ws.Range("R13").Value = 185.545940907725
$ws.Range("S13").Value = 0.01250218376717629
$ws.Range("T13").Value = 0.01250218376717629
